# Generate Report for Handback
# The e090e645-09b3-4542-97c9-c769b01b38fb file has completed its handback
# round-trip: update Status to "Handed back: in sync with en-US" on every
# sheet, and fill in the Latest Target File / Latest Handback File /
# Latest Handback DateTime columns on the per-language detail sheets.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # OLE BGR for RGB(100,149,237) == ARGB FF6495ED used by the "HyperLink" style

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns for the e090e645 row
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = $statusDone
$wsOverview.Range("F6").Value = $statusDone

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = $statusDone
$wsZhCn.Range("I6").Value = "e090e645-09b3-4542-97c9-c769b01b38fb.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/e090e645-09b3-4542-97c9-c769b01b38fb.md", $null, $null, "e090e645-09b3-4542-97c9-c769b01b38fb.md")
$wsZhCn.Range("I6").Font.Underline = $true
$wsZhCn.Range("I6").Font.Color = $hyperlinkColor
$wsZhCn.Range("J6").Value = "e090e645-09b3-4542-97c9-c769b01b38fb.22712eecc61c5fc8df56fcf01097a08c1630d156.zh-cn.xlf"
$wsZhCn.Range("K6").Value = "2016-10-10 09:47:59"

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = $statusDone
$wsDeDe.Range("I6").Value = "e090e645-09b3-4542-97c9-c769b01b38fb.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/e090e645-09b3-4542-97c9-c769b01b38fb.md", $null, $null, "e090e645-09b3-4542-97c9-c769b01b38fb.md")
$wsDeDe.Range("I6").Font.Underline = $true
$wsDeDe.Range("I6").Font.Color = $hyperlinkColor
$wsDeDe.Range("J6").Value = "e090e645-09b3-4542-97c9-c769b01b38fb.22712eecc61c5fc8df56fcf01097a08c1630d156.de-de.xlf"
$wsDeDe.Range("K6").Value = "2016-10-10 09:48:15"
